$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Project Total Costs")
$ws.Range("B2").Value = [double]"881252.825989"
$ws.Range("B3").Value = [double]"181252.968612"
$ws.Range("B4").Value = [double]"707961.514072"
$ws.Range("B5").Value = [double]"699999.857377"
$ws.Range("B7").Value = [double]"0.6804938176950782"

$ws = $wb.Worksheets.Item("Components Capacity and Cost")
$ws.Range("B2").Value = [double]"532.8082786489999"
$ws.Range("B3").Value = [double]"28.0624080285"
$ws.Range("B4").Value = [double]"5.53638950717e-05"
$ws.Range("B5").Value = [double]"171.831087874"
$ws.Range("B6").Value = [double]"21.0028877056"
$ws.Range("B7").Value = [double]"287716.47047046"
$ws.Range("B8").Value = [double]"5612.4816057"
$ws.Range("B9").Value = [double]"0.06643667408604"
$ws.Range("B10").Value = [double]"343662.175748"
$ws.Range("B11").Value = [double]"63008.66311679999"
$ws.Range("B12").Value = [double]"699999.857377634"

$ws = $wb.Worksheets.Item("Yearly Costs Info")
$ws.Range("B2").Value = [double]"5754.329409409201"
$ws.Range("C2").Value = [double]"561.2780570733387"
$ws.Range("D2").Value = [double]"8133.416777295999"
$ws.Range("E2").Value = [double]"14449.02424377854"
$ws.Range("F2").Value = [double]"0.04125488801850326"
$ws.Range("G2").Value = [double]"8893.018399208784"
$ws.Range("B3").Value = [double]"5754.329409409201"
$ws.Range("C3").Value = [double]"561.2780570733387"
$ws.Range("D3").Value = [double]"8133.416777295999"
$ws.Range("E3").Value = [double]"14449.02424377854"
$ws.Range("F3").Value = [double]"0.04125193280751012"
$ws.Range("G3").Value = [double]"8911.786229106523"
$ws.Range("B4").Value = [double]"5754.329409409201"
$ws.Range("C4").Value = [double]"561.2780570733387"
$ws.Range("D4").Value = [double]"8133.416777295999"
$ws.Range("E4").Value = [double]"14449.02424377854"
$ws.Range("F4").Value = [double]"76.47893717906528"
$ws.Range("G4").Value = [double]"11024.66356877324"
$ws.Range("B5").Value = [double]"5754.329409409201"
$ws.Range("C5").Value = [double]"561.2780570733387"
$ws.Range("D5").Value = [double]"8133.416777295999"
$ws.Range("E5").Value = [double]"14449.02424377854"
$ws.Range("F5").Value = [double]"81.92750187927129"
$ws.Range("G5").Value = [double]"11024.87703796586"
$ws.Range("B6").Value = [double]"5754.329409409201"
$ws.Range("C6").Value = [double]"561.2780570733387"
$ws.Range("D6").Value = [double]"8133.416777295999"
$ws.Range("E6").Value = [double]"14449.02424377854"
$ws.Range("F6").Value = [double]"81.92750187927129"
$ws.Range("G6").Value = [double]"11024.87703796586"
$ws.Range("B7").Value = [double]"5754.329409409201"
$ws.Range("C7").Value = [double]"561.2780570733387"
$ws.Range("D7").Value = [double]"8133.416777295999"
$ws.Range("E7").Value = [double]"14449.02424377854"
$ws.Range("F7").Value = [double]"81.92750187927129"
$ws.Range("G7").Value = [double]"11024.87703796586"
$ws.Range("B8").Value = [double]"5754.329409409201"
$ws.Range("C8").Value = [double]"561.2780570733387"
$ws.Range("D8").Value = [double]"8133.416777295999"
$ws.Range("E8").Value = [double]"14449.02424377854"
$ws.Range("F8").Value = [double]"100.5919460951342"
$ws.Range("G8").Value = [double]"11021.40227621262"
$ws.Range("B9").Value = [double]"5754.329409409201"
$ws.Range("C9").Value = [double]"561.2780570733387"
$ws.Range("D9").Value = [double]"8133.416777295999"
$ws.Range("E9").Value = [double]"14449.02424377854"
$ws.Range("F9").Value = [double]"3717.077487800779"
$ws.Range("G9").Value = [double]"14751.75318293292"
$ws.Range("B10").Value = [double]"5754.329409409201"
$ws.Range("C10").Value = [double]"561.2780570733387"
$ws.Range("D10").Value = [double]"8133.416777295999"
$ws.Range("E10").Value = [double]"14449.02424377854"
$ws.Range("F10").Value = [double]"3739.768033320368"
$ws.Range("G10").Value = [double]"14747.52888169429"
$ws.Range("B11").Value = [double]"5754.329409409201"
$ws.Range("C11").Value = [double]"561.2780570733387"
$ws.Range("D11").Value = [double]"8133.416777295999"
$ws.Range("E11").Value = [double]"14449.02424377854"
$ws.Range("F11").Value = [double]"3739.768033320368"
$ws.Range("G11").Value = [double]"14747.52888169429"
$ws.Range("B12").Value = [double]"5754.329409409201"
$ws.Range("C12").Value = [double]"561.2780570733387"
$ws.Range("D12").Value = [double]"8133.416777295999"
$ws.Range("E12").Value = [double]"14449.02424377854"
$ws.Range("F12").Value = [double]"3739.768033320368"
$ws.Range("G12").Value = [double]"14747.52888169429"
$ws.Range("B13").Value = [double]"5754.329409409201"
$ws.Range("C13").Value = [double]"561.2780570733387"
$ws.Range("D13").Value = [double]"8133.416777295999"
$ws.Range("E13").Value = [double]"14449.02424377854"
$ws.Range("F13").Value = [double]"3739.768033320368"
$ws.Range("G13").Value = [double]"14747.52888169429"
$ws.Range("B14").Value = [double]"5754.329409409201"
$ws.Range("C14").Value = [double]"561.2780570733387"
$ws.Range("D14").Value = [double]"8133.416777295999"
$ws.Range("E14").Value = [double]"14449.02424377854"
$ws.Range("F14").Value = [double]"3739.768033320368"
$ws.Range("G14").Value = [double]"14747.52888169429"
$ws.Range("B15").Value = [double]"5754.329409409201"
$ws.Range("C15").Value = [double]"561.2780570733387"
$ws.Range("D15").Value = [double]"8133.416777295999"
$ws.Range("E15").Value = [double]"14449.02424377854"
$ws.Range("F15").Value = [double]"3765.578516434562"
$ws.Range("G15").Value = [double]"14742.72374135239"
$ws.Range("B16").Value = [double]"5754.329409409201"
$ws.Range("C16").Value = [double]"561.2780570733387"
$ws.Range("D16").Value = [double]"8133.416777295999"
$ws.Range("E16").Value = [double]"14449.02424377854"
$ws.Range("F16").Value = [double]"21976.74275005401"
$ws.Range("G16").Value = [double]"14029.6055309928"
$ws.Range("B17").Value = [double]"5754.329409409201"
$ws.Range("C17").Value = [double]"561.2780570733387"
$ws.Range("D17").Value = [double]"8133.416777295999"
$ws.Range("E17").Value = [double]"14449.02424377854"
$ws.Range("F17").Value = [double]"22011.83316269384"
$ws.Range("G17").Value = [double]"14033.38187692725"
$ws.Range("B18").Value = [double]"5754.329409409201"
$ws.Range("C18").Value = [double]"561.2780570733387"
$ws.Range("D18").Value = [double]"8133.416777295999"
$ws.Range("E18").Value = [double]"14449.02424377854"
$ws.Range("F18").Value = [double]"22011.83316269384"
$ws.Range("G18").Value = [double]"14033.38187692725"
$ws.Range("B19").Value = [double]"5754.329409409201"
$ws.Range("C19").Value = [double]"561.2780570733387"
$ws.Range("D19").Value = [double]"8133.416777295999"
$ws.Range("E19").Value = [double]"14449.02424377854"
$ws.Range("F19").Value = [double]"22011.83316269384"
$ws.Range("G19").Value = [double]"14033.38187692725"
$ws.Range("B20").Value = [double]"5754.329409409201"
$ws.Range("C20").Value = [double]"561.2780570733387"
$ws.Range("D20").Value = [double]"8133.416777295999"
$ws.Range("E20").Value = [double]"14449.02424377854"
$ws.Range("F20").Value = [double]"22011.83316269408"
$ws.Range("G20").Value = [double]"14033.38187692725"
$ws.Range("B21").Value = [double]"5754.329409409201"
$ws.Range("C21").Value = [double]"561.2780570733387"
$ws.Range("D21").Value = [double]"8133.416777295999"
$ws.Range("E21").Value = [double]"14449.02424377854"
$ws.Range("F21").Value = [double]"22000.49476444961"
$ws.Range("G21").Value = [double]"14031.27100656734"

$ws = $wb.Worksheets.Item("Yearly Energy Averages")
$ws.Range("B2").Value = [double]"5.633228088194637e-05"
$ws.Range("C2").Value = [double]"5.969367661581906e-05"
$ws.Range("D2").Value = [double]"106.0452606407318"
$ws.Range("E2").Value = [double]"55.23491783702503"
$ws.Range("B3").Value = [double]"5.632828836249129e-05"
$ws.Range("C3").Value = [double]"5.968934895265128e-05"
$ws.Range("D3").Value = [double]"106.3036448558265"
$ws.Range("E3").Value = [double]"55.12584575131432"
$ws.Range("B4").Value = [double]"0.1540634758446386"
$ws.Range("C4").Value = [double]"4.822266525112051e-05"
$ws.Range("D4").Value = [double]"105.9368695933263"
$ws.Range("E4").Value = [double]"43.08082674833905"
$ws.Range("B5").Value = [double]"0.165042130587793"
$ws.Range("C5").Value = [double]"4.832436640052904e-05"
$ws.Range("D5").Value = [double]"105.9501588963825"
$ws.Range("E5").Value = [double]"43.07115987036232"
$ws.Range("B6").Value = [double]"0.165042130587793"
$ws.Range("C6").Value = [double]"4.832436640052904e-05"
$ws.Range("D6").Value = [double]"105.9501588963825"
$ws.Range("E6").Value = [double]"43.07115987036232"
$ws.Range("B7").Value = [double]"0.165042130587793"
$ws.Range("C7").Value = [double]"4.832436640052904e-05"
$ws.Range("D7").Value = [double]"105.9501588963825"
$ws.Range("E7").Value = [double]"43.07115987036232"
$ws.Range("B8").Value = [double]"0.2026505834167316"
$ws.Range("C8").Value = [double]"4.831806941566382e-05"
$ws.Range("D8").Value = [double]"105.9501588934641"
$ws.Range("E8").Value = [double]"43.06247149041093"
$ws.Range("B9").Value = [double]"5.107667908166076"
$ws.Range("C9").Value = [double]"6.415729257910358e-05"
$ws.Range("D9").Value = [double]"100.4395554794899"
$ws.Range("E9").Value = [double]"20.10728594638768"
$ws.Range("B10").Value = [double]"5.13884750610644"
$ws.Range("C10").Value = [double]"6.415550227133429e-05"
$ws.Range("D10").Value = [double]"100.4395554757131"
$ws.Range("E10").Value = [double]"20.10254154646164"
$ws.Range("B11").Value = [double]"5.13884750610644"
$ws.Range("C11").Value = [double]"6.415550227133429e-05"
$ws.Range("D11").Value = [double]"100.4395554757131"
$ws.Range("E11").Value = [double]"20.10254154646164"
$ws.Range("B12").Value = [double]"5.13884750610644"
$ws.Range("C12").Value = [double]"6.415550227133429e-05"
$ws.Range("D12").Value = [double]"100.4395554757131"
$ws.Range("E12").Value = [double]"20.10254154646164"
$ws.Range("B13").Value = [double]"5.13884750610644"
$ws.Range("C13").Value = [double]"6.415550227133429e-05"
$ws.Range("D13").Value = [double]"100.4395554757131"
$ws.Range("E13").Value = [double]"20.10254154646164"
$ws.Range("B14").Value = [double]"5.13884750610644"
$ws.Range("C14").Value = [double]"6.415550227133429e-05"
$ws.Range("D14").Value = [double]"100.4395554757131"
$ws.Range("E14").Value = [double]"20.10254154646164"
$ws.Range("B15").Value = [double]"5.174314281794326"
$ws.Range("C15").Value = [double]"6.415335760299953e-05"
$ws.Range("D15").Value = [double]"100.4395554741389"
$ws.Range("E15").Value = [double]"20.09714751321988"
$ws.Range("B16").Value = [double]"22.14823420263563"
$ws.Range("C16").Value = [double]"6.693151502430933e-05"
$ws.Range("D16").Value = [double]"81.70955748336661"
$ws.Range("E16").Value = [double]"9.919843905827786"
$ws.Range("B17").Value = [double]"22.18359837544556"
$ws.Range("C17").Value = [double]"6.711254800830729e-05"
$ws.Range("D17").Value = [double]"81.70955753668488"
$ws.Range("E17").Value = [double]"9.916802084210302"
$ws.Range("B18").Value = [double]"22.18359837544556"
$ws.Range("C18").Value = [double]"6.711254800830729e-05"
$ws.Range("D18").Value = [double]"81.70955753668488"
$ws.Range("E18").Value = [double]"9.916802084210302"
$ws.Range("B19").Value = [double]"22.18359837544556"
$ws.Range("C19").Value = [double]"6.711254800830729e-05"
$ws.Range("D19").Value = [double]"81.70955753668488"
$ws.Range("E19").Value = [double]"9.916802084210302"
$ws.Range("B20").Value = [double]"22.18359837544581"
$ws.Range("C20").Value = [double]"6.711254800830851e-05"
$ws.Range("D20").Value = [double]"81.70955753668488"
$ws.Range("E20").Value = [double]"9.916802084210282"
$ws.Range("B21").Value = [double]"22.17217147041952"
$ws.Range("C21").Value = [double]"6.711388217510287e-05"
$ws.Range("D21").Value = [double]"81.70955753701425"
$ws.Range("E21").Value = [double]"9.917784736101492"
